$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cA = $ws.Cells.Item(2, 1)
$cA.Value = "'2024-12-22"
$cA.Style = "Normal"
$ws.Cells.Item(2, 2).Value = 0.1516480774644057
$ws.Cells.Item(2, 3).Value = 0.09841994491865877
$ws.Cells.Item(2, 4).Value = 1.388418552533649
$ws.Cells.Item(2, 5).Value = "{'2330.TW': 0.10046, '2449.TW': 0.08815, '6669.TW': 0.08016, '2474.TW': 0.07576, '3702.TW': 0.07042}"

$cA = $ws.Cells.Item(3, 1)
$cA.Value = "'2025-01-05"
$cA.Style = "Normal"
$ws.Cells.Item(3, 2).Value = 0.1516480774644057
$ws.Cells.Item(3, 3).Value = 0.1003019146413577
$ws.Cells.Item(3, 4).Value = 1.362367587428499
$ws.Cells.Item(3, 5).Value = "{'2330.TW': 0.09495, '2449.TW': 0.08113, '6669.TW': 0.08079, '3702.TW': 0.07896, '2474.TW': 0.07139}"

$cA = $ws.Cells.Item(4, 1)
$cA.Value = "'2025-01-19"
$cA.Style = "Normal"
$ws.Cells.Item(4, 2).Value = 0.1516480774644057
$ws.Cells.Item(4, 3).Value = 0.1009193956496041
$ws.Cells.Item(4, 4).Value = 1.354031864586792
$ws.Cells.Item(4, 5).Value = "{'2330.TW': 0.10031, '6669.TW': 0.08423, '2449.TW': 0.08299, '3702.TW': 0.08122, '2474.TW': 0.06991}"

$cA = $ws.Cells.Item(5, 1)
$cA.Value = "'2025-02-02"
$cA.Style = "Normal"
$ws.Cells.Item(5, 2).Value = 0.1516480774644057
$ws.Cells.Item(5, 3).Value = 0.0994157972633072
$ws.Cells.Item(5, 4).Value = 1.374510703791743
$ws.Cells.Item(5, 5).Value = "{'2330.TW': 0.09815, '2449.TW': 0.0841, '3702.TW': 0.08041, '6669.TW': 0.07918, '2474.TW': 0.06798}"

$cA = $ws.Cells.Item(6, 1)
$cA.Value = "'2025-02-16"
$cA.Style = "Normal"
$ws.Cells.Item(6, 2).Value = 0.1516480774644057
$ws.Cells.Item(6, 3).Value = 0.09756232121158734
$ws.Cells.Item(6, 4).Value = 1.400623476024637
$ws.Cells.Item(6, 5).Value = "{'2330.TW': 0.10165, '2449.TW': 0.08159, '6669.TW': 0.07751, '3702.TW': 0.07528, '2474.TW': 0.06702}"

$cA = $ws.Cells.Item(7, 1)
$cA.Value = "'2025-03-02"
$cA.Style = "Normal"
$ws.Cells.Item(7, 2).Value = 0.1516480774644056
$ws.Cells.Item(7, 3).Value = 0.093869136025188
$ws.Cells.Item(7, 4).Value = 1.455729574710677
$ws.Cells.Item(7, 5).Value = "{'2330.TW': 0.09895, '2449.TW': 0.07856, '6669.TW': 0.07496, '3702.TW': 0.07415, '2474.TW': 0.06688}"

$cA = $ws.Cells.Item(8, 1)
$cA.Value = "'2025-03-16"
$cA.Style = "Normal"
$ws.Cells.Item(8, 2).Value = 0.1516480774644057
$ws.Cells.Item(8, 3).Value = 0.09249499340045338
$ws.Cells.Item(8, 4).Value = 1.477356475639642
$ws.Cells.Item(8, 5).Value = "{'2330.TW': 0.10076, '2449.TW': 0.07817, '6669.TW': 0.07403, '3702.TW': 0.07351, '2474.TW': 0.07133}"

$cA = $ws.Cells.Item(9, 1)
$cA.Value = "'2025-03-30"
$cA.Style = "Normal"
$ws.Cells.Item(9, 2).Value = 0.1516480774644056
$ws.Cells.Item(9, 3).Value = 0.09185067239047918
$ws.Cells.Item(9, 4).Value = 1.487719947040583
$ws.Cells.Item(9, 5).Value = "{'2330.TW': 0.11101, '2449.TW': 0.07501, '3702.TW': 0.07415, '6669.TW': 0.06942, '2474.TW': 0.06705}"

$cA = $ws.Cells.Item(10, 1)
$cA.Value = "'2025-04-13"
$cA.Style = "Normal"
$ws.Cells.Item(10, 2).Value = 0.1516480774644056
$ws.Cells.Item(10, 3).Value = 0.1743576838845118
$ws.Cells.Item(10, 4).Value = 0.7837227153975975
$ws.Cells.Item(10, 5).Value = "{'2449.TW': 0.23143, '2330.TW': 0.157, '1504.TW': 0.14353, '2474.TW': 0.10503, '3702.TW': 0.09154}"

$cA = $ws.Cells.Item(11, 1)
$cA.Value = "'2025-04-27"
$cA.Style = "Normal"
$ws.Cells.Item(11, 2).Value = 0.1516480774644056
$ws.Cells.Item(11, 3).Value = 0.1669854087581121
$ws.Cells.Item(11, 4).Value = 0.8183234599997186
$ws.Cells.Item(11, 5).Value = "{'2449.TW': 0.21053, '1504.TW': 0.1449, '2330.TW': 0.14223, '2474.TW': 0.10221, '3017.TW': 0.09631}"

$cA = $ws.Cells.Item(12, 1)
$cA.Value = "'2025-05-11"
$cA.Style = "Normal"
$ws.Cells.Item(12, 2).Value = 0.1516480774644056
$ws.Cells.Item(12, 3).Value = 0.1591811668596504
$ws.Cells.Item(12, 4).Value = 0.8584437478391388
$ws.Cells.Item(12, 5).Value = "{'2449.TW': 0.20328, '1504.TW': 0.14811, '2330.TW': 0.11889, '3017.TW': 0.10258, '2474.TW': 0.10179}"

$cA = $ws.Cells.Item(13, 1)
$cA.Value = "'2025-05-25"
$cA.Style = "Normal"
$ws.Cells.Item(13, 2).Value = 0.1516480774644056
$ws.Cells.Item(13, 3).Value = 0.1526547830051917
$ws.Cells.Item(13, 4).Value = 0.8951444217752307
$ws.Cells.Item(13, 5).Value = "{'2449.TW': 0.19341, '1504.TW': 0.14442, '2330.TW': 0.11183, '3017.TW': 0.10496, '2474.TW': 0.0909}"

$cA = $ws.Cells.Item(14, 1)
$cA.Value = "'2025-06-08"
$cA.Style = "Normal"
$ws.Cells.Item(14, 2).Value = 0.1516480774644056
$ws.Cells.Item(14, 3).Value = 0.1484084686107517
$ws.Cells.Item(14, 4).Value = 0.9207566033364883
$ws.Cells.Item(14, 5).Value = "{'2449.TW': 0.18149, '1504.TW': 0.14117, '3017.TW': 0.10014, '2330.TW': 0.09647, '3702.TW': 0.08495}"

$cA = $ws.Cells.Item(15, 1)
$cA.Value = "'2025-06-22"
$cA.Style = "Normal"
$ws.Cells.Item(15, 2).Value = 0.1516480774644056
$ws.Cells.Item(15, 3).Value = 0.1430413222794237
$ws.Cells.Item(15, 4).Value = 0.9553049097062372
$ws.Cells.Item(15, 5).Value = "{'2449.TW': 0.17217, '1504.TW': 0.13281, '2330.TW': 0.09756, '3017.TW': 0.09327, '3702.TW': 0.08708}"

$cA = $ws.Cells.Item(16, 1)
$cA.Value = "'2025-07-06"
$cA.Style = "Normal"
$ws.Cells.Item(16, 2).Value = 0.1516480774644056
$ws.Cells.Item(16, 3).Value = 0.1399883199739695
$ws.Cells.Item(16, 4).Value = 0.9761391342493072
$ws.Cells.Item(16, 5).Value = "{'2449.TW': 0.1534, '1504.TW': 0.12779, '3017.TW': 0.08636, '2330.TW': 0.07912, '3702.TW': 0.07708}"

$ws.Range("A17:E17").ClearContents()